# M4 and M5 slide updates
#
# Slide 9 ("Creating an Azure Mobile Service / Adding data / Creating API
# method") has a speaker-notes paragraph that used to end in two adjacent
# runs with identical-looking formatting:
#   run4: ", show the PowerShell functions for simple get "
#   run5: "and post"
# The author selected across that run boundary and retyped it as a single
# run:
#   ", show the PowerShell functions for simple get and post"
# which merges run4+run5 (PowerPoint coalesces the now-identically-typed
# text into one <a:r>) and refreshes the paragraph's trailing endParaRPr to
# dirty="0". The other runs in the paragraph are untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$notesPage = $s.NotesPage
$notesShape = $notesPage.Shapes.Item(2)
$tr = $notesShape.TextFrame.TextRange

# Reconstruct the paragraph text with only the run4/run5 boundary removed
# (no other characters change).
$run1 = "we're going to do a lot here " + [char]0x2013 + " create the mobile service, create a knots table, create a steps table,"
$run2 = " add sample data, create a knots "
$run3 = "api"
$run4 = ", show the PowerShell functions for simple get "
$run5 = "and post"

$mergedTail = $run4 + $run5
$newText = $run1 + $run2 + $run3 + $mergedTail

$tr.Text = $newText
